# ------------------------------------------------------------------
# Edit script: updates resume body text (work experience blurb),
# applies Times New Roman font throughout body/header/footer, and
# bumps the "Last updated on" date in the header.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Text edits inside the body paragraph -----------------------

# "23 Months" -> "2 Years"
$d.Content.Find.Execute(
    "23 Months of total work experience in Infosys Ltd.", $false, $true,
    $false, $false, $false, $true, 1, $false,
    "2 Years of total work experience in Infosys Ltd.", 2) | Out-Null

# "Systems Engineer" -> "Senior Systems Engineer"
$d.Content.Find.Execute(
    "Currently working as Systems Engineer", $false, $true,
    $false, $false, $false, $true, 1, $false,
    "Currently working as Senior Systems Engineer", 2) | Out-Null

# Mention new PLM tool before Selerant's DevEX PLM (note curly apostrophe
# in the replacement, matching the author's smart-quote autocorrect).
$d.Content.Find.Execute(
    "Meridian, Selerant's DevEX PLM.", $false, $true,
    $false, $false, $false, $true, 1, $false,
    "Meridian, Oracle Agile PLM, and Selerant" + [char]0x2019 + "s DevEX PLM.", 2) | Out-Null

# --- 2. Restore the "_GoBack" last-edit bookmark at its new spot ---
# (Word always leaves this bookmark at wherever text was last typed;
# after the edits above that is inside "...data structures, databas|e,
# algorithms...".)

$rng = $d.Content
$null = $rng.Find.Execute("structures, databas", $false, $false, $false,
    $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 3. Apply Times New Roman (ascii/hAnsi/cs) to the whole body ---

$body = $d.Paragraphs.Item(1).Range
$body.Font.Name = "Times New Roman"
$body.Font.NameBi = "Times New Roman"

# --- 4. Footer: font + size ------------------------------------------

$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Font.Name = "Times New Roman"
$footer.Range.Font.NameBi = "Times New Roman"
$footer.Range.Font.Size = 17
$footer.Range.Font.SizeBi = 17

# --- 5. Header: font + size + updated date ---------------------------

$header = $d.Sections.Item(1).Headers.Item(1)
$header.Range.Font.Name = "Times New Roman"
$header.Range.Font.NameBi = "Times New Roman"
$header.Range.Font.Size = 17
$header.Range.Font.SizeBi = 17

$header.Range.Find.Execute(
    "May 23, 2017", $false, $true,
    $false, $false, $false, $true, 1, $false,
    "July 21, 2017", 2) | Out-Null
